$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 202.28572  # ALC!H2: was 187
$ws.Cells.Item(2, 9).Value = 202.28572  # ALC!I2: was 187
$ws.Cells.Item(2, 11).Value = 202.28572  # ALC!K2: was 187
$ws.Cells.Item(2, 13).Value = -89.28572  # ALC!M2: was -74

$ws.Cells.Item(12, 8).Value = 4166.5835  # ALC!H12: was 4454.4546
$ws.Cells.Item(12, 9).Value = 6333.3335  # ALC!I12: was 7400
$ws.Cells.Item(12, 11).Value = 6333.3335  # ALC!K12: was 7400
$ws.Cells.Item(12, 13).Value = -6163.3335  # ALC!M12: was -7230

$ws.Cells.Item(31, 8).Value = 542.7143  # ALC!H31: was 542.8570999999999
$ws.Cells.Item(31, 9).Value = 466.5  # ALC!I31: was 466.66666
$ws.Cells.Item(31, 11).Value = 1399.5  # ALC!K31: was 1399.99998
$ws.Cells.Item(31, 13).Value = -1169.5  # ALC!M31: was -1169.99998

$ws.Cells.Item(38, 8).Value = 517.3  # ALC!H38: was 565.7273
$ws.Cells.Item(38, 9).Value = 517.3  # ALC!I38: was 565.7273
$ws.Cells.Item(38, 11).Value = 1551.9  # ALC!K38: was 1697.1819
$ws.Cells.Item(38, 13).Value = -1179.9  # ALC!M38: was -1325.1819

$ws.Cells.Item(70, 8).Value = 3450.4443  # ALC!H70: was 3168.5454
$ws.Cells.Item(70, 10).Value = 2816.6  # ALC!J70: was 2554.7144
$ws.Cells.Item(70, 12).Value = 8449.799999999999  # ALC!L70: was 7664.1432
$ws.Cells.Item(70, 14).Value = -8989.799999999999  # ALC!N70: was -8204.143199999999

$ws.Cells.Item(73, 8).Value = 3450.4443  # ALC!H73: was 3168.5454
$ws.Cells.Item(73, 10).Value = 2816.6  # ALC!J73: was 2554.7144
$ws.Cells.Item(73, 12).Value = 8449.799999999999  # ALC!L73: was 7664.1432
$ws.Cells.Item(73, 14).Value = -10321.8  # ALC!N73: was -9536.143199999999

$ws.Cells.Item(92, 8).Value = 351.81818  # ALC!H92: was 364.61905
$ws.Cells.Item(92, 9).Value = 351.81818  # ALC!I92: was 364.61905
$ws.Cells.Item(92, 11).Value = 351.81818  # ALC!K92: was 364.61905
$ws.Cells.Item(92, 13).Value = 896.18182  # ALC!M92: was 883.38095

$ws.Cells.Item(127, 8).Value = 931.2  # ALC!H127: was 894.5
$ws.Cells.Item(127, 9).Value = 931.2  # ALC!I127: was 879.4286
$ws.Cells.Item(127, 10).Value = 0  # ALC!J127: was 1000
$ws.Cells.Item(127, 11).Value = 2793.6  # ALC!K127: was 2638.2858
$ws.Cells.Item(127, 12).Value = 0  # ALC!L127: was 3000
$ws.Cells.Item(127, 13).Value = 2166.4  # ALC!M127: was 2321.7142
$ws.Cells.Item(127, 14).ClearContents()  # ALC!N127: was -12920

$ws.Cells.Item(138, 8).Value = 2117.2917  # ALC!H138: was 2180.1333
$ws.Cells.Item(138, 9).Value = 1717.6923  # ALC!I138: was 1751
$ws.Cells.Item(138, 10).Value = 2589.5454  # ALC!J138: was 2555.625
$ws.Cells.Item(138, 11).Value = 5153.0769  # ALC!K138: was 5253
$ws.Cells.Item(138, 12).Value = 7768.6362  # ALC!L138: was 7666.875
$ws.Cells.Item(138, 13).Value = -13.07690000000002  # ALC!M138: was -113
$ws.Cells.Item(138, 14).Value = -18048.6362  # ALC!N138: was -17946.875

$ws.Cells.Item(140, 8).Value = 74994.5  # ALC!H140: was 110256.336
$ws.Cells.Item(140, 10).Value = 74994.5  # ALC!J140: was 110256.336
$ws.Cells.Item(140, 12).Value = 74994.5  # ALC!L140: was 110256.336
$ws.Cells.Item(140, 14).Value = -85354.5  # ALC!N140: was -120616.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6451.9697  # ARM!H32: was 6653.7188
$ws.Cells.Item(32, 9).Value = 2623.9375  # ARM!I32: was 2708.7097
$ws.Cells.Item(32, 11).Value = 2623.9375  # ARM!K32: was 2708.7097
$ws.Cells.Item(32, 13).Value = -2336.9375  # ARM!M32: was -2421.7097

$ws.Cells.Item(45, 8).Value = 14661.889  # ARM!H45: was 12641.091
$ws.Cells.Item(45, 9).Value = 35469  # ARM!I45: was 27476.75
$ws.Cells.Item(45, 10).Value = 4258.3335  # ARM!J45: was 4163.5713
$ws.Cells.Item(45, 11).Value = 35469  # ARM!K45: was 27476.75
$ws.Cells.Item(45, 12).Value = 4258.3335  # ARM!L45: was 4163.5713
$ws.Cells.Item(45, 13).Value = -35092  # ARM!M45: was -27099.75
$ws.Cells.Item(45, 14).Value = -5012.3335  # ARM!N45: was -4917.5713

$ws.Cells.Item(61, 8).Value = 1907.75  # ARM!H61: was 2008.8572
$ws.Cells.Item(61, 9).Value = 1685.2  # ARM!I61: was 1806.5
$ws.Cells.Item(61, 11).Value = 1685.2  # ARM!K61: was 1806.5
$ws.Cells.Item(61, 13).Value = -1473.2  # ARM!M61: was -1594.5

$ws.Cells.Item(97, 8).Value = 1678.421  # ARM!H97: was 2140.9412
$ws.Cells.Item(97, 9).Value = 1730  # ARM!I97: was 2227.875
$ws.Cells.Item(97, 11).Value = 1730  # ARM!K97: was 2227.875
$ws.Cells.Item(97, 13).Value = -1234  # ARM!M97: was -1731.875

$ws.Cells.Item(110, 8).Value = 1799.75  # ARM!H110: was 1875
$ws.Cells.Item(110, 9).Value = 1799.75  # ARM!I110: was 1875
$ws.Cells.Item(110, 11).Value = 1799.75  # ARM!K110: was 1875
$ws.Cells.Item(110, 13).Value = 245.25  # ARM!M110: was 170

$ws.Cells.Item(136, 8).Value = 1907.75  # ARM!H136: was 2008.8572
$ws.Cells.Item(136, 9).Value = 1685.2  # ARM!I136: was 1806.5
$ws.Cells.Item(136, 11).Value = 5055.6  # ARM!K136: was 5419.5
$ws.Cells.Item(136, 13).Value = -2505.6  # ARM!M136: was -2869.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 12352.857  # BSM!H20: was 12921.143
$ws.Cells.Item(20, 9).Value = 11705.1875  # BSM!I20: was 13094.6
$ws.Cells.Item(20, 10).Value = 14425.4  # BSM!J20: was 12487.5
$ws.Cells.Item(20, 11).Value = 11705.1875  # BSM!K20: was 13094.6
$ws.Cells.Item(20, 12).Value = 14425.4  # BSM!L20: was 12487.5
$ws.Cells.Item(20, 13).Value = -11458.1875  # BSM!M20: was -12847.6
$ws.Cells.Item(20, 14).Value = -14919.4  # BSM!N20: was -12981.5

$ws.Cells.Item(86, 8).Value = 2021.6842  # BSM!H86: was 2465.0588
$ws.Cells.Item(86, 9).Value = 1755.5454  # BSM!I86: was 2336.818
$ws.Cells.Item(86, 10).Value = 2387.625  # BSM!J86: was 2700.1667
$ws.Cells.Item(86, 11).Value = 1755.5454  # BSM!K86: was 2336.818
$ws.Cells.Item(86, 12).Value = 2387.625  # BSM!L86: was 2700.1667
$ws.Cells.Item(86, 13).Value = -632.5454  # BSM!M86: was -1213.818
$ws.Cells.Item(86, 14).Value = -4633.625  # BSM!N86: was -4946.1667

$ws.Cells.Item(89, 8).Value = 2021.6842  # BSM!H89: was 2465.0588
$ws.Cells.Item(89, 9).Value = 1755.5454  # BSM!I89: was 2336.818
$ws.Cells.Item(89, 10).Value = 2387.625  # BSM!J89: was 2700.1667
$ws.Cells.Item(89, 11).Value = 8777.726999999999  # BSM!K89: was 11684.09
$ws.Cells.Item(89, 12).Value = 11938.125  # BSM!L89: was 13500.8335
$ws.Cells.Item(89, 13).Value = -3161.726999999999  # BSM!M89: was -6068.09
$ws.Cells.Item(89, 14).Value = -23170.125  # BSM!N89: was -24732.8335

$ws.Cells.Item(94, 8).Value = 1986.1666  # BSM!H94: was 768.6667
$ws.Cells.Item(94, 9).Value = 2316.75  # BSM!I94: was 829.8125
$ws.Cells.Item(94, 10).Value = 1325  # BSM!J94: was 279.5
$ws.Cells.Item(94, 11).Value = 2316.75  # BSM!K94: was 829.8125
$ws.Cells.Item(94, 12).Value = 1325  # BSM!L94: was 279.5
$ws.Cells.Item(94, 13).Value = -1865.75  # BSM!M94: was -378.8125
$ws.Cells.Item(94, 14).Value = -2227  # BSM!N94: was -1181.5

$ws.Cells.Item(134, 8).Value = 1722  # BSM!H134: was 1852.7142
$ws.Cells.Item(134, 9).Value = 1429.3334  # BSM!I134: was 1553.8
$ws.Cells.Item(134, 11).Value = 4288.0002  # BSM!K134: was 4661.4
$ws.Cells.Item(134, 13).Value = -1753.0002  # BSM!M134: was -2126.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(10, 8).Value = 1171.5555  # CRP!H10: was 1205.25
$ws.Cells.Item(10, 9).Value = 770.5714  # CRP!I10: was 732.3333
$ws.Cells.Item(10, 10).Value = 2575  # CRP!J10: was 2624
$ws.Cells.Item(10, 11).Value = 770.5714  # CRP!K10: was 732.3333
$ws.Cells.Item(10, 12).Value = 2575  # CRP!L10: was 2624
$ws.Cells.Item(10, 13).Value = -631.5714  # CRP!M10: was -593.3333
$ws.Cells.Item(10, 14).Value = -2853  # CRP!N10: was -2902

$ws.Cells.Item(11, 8).Value = 500377.25  # CRP!H11: was 333635.16
$ws.Cells.Item(11, 9).Value = 500  # CRP!I11: was 252.5
$ws.Cells.Item(11, 10).Value = 667003  # CRP!J11: was 500326.5
$ws.Cells.Item(11, 11).Value = 500  # CRP!K11: was 252.5
$ws.Cells.Item(11, 12).Value = 667003  # CRP!L11: was 500326.5
$ws.Cells.Item(11, 13).Value = -360  # CRP!M11: was -112.5
$ws.Cells.Item(11, 14).Value = -667283  # CRP!N11: was -500606.5

$ws.Cells.Item(29, 8).Value = 0  # CRP!H29: was 649.75
$ws.Cells.Item(29, 9).Value = 0  # CRP!I29: was 500
$ws.Cells.Item(29, 10).Value = 0  # CRP!J29: was 699.6667
$ws.Cells.Item(29, 11).Value = 0  # CRP!K29: was 500
$ws.Cells.Item(29, 12).Value = 0  # CRP!L29: was 699.6667
$ws.Cells.Item(29, 13).ClearContents()  # CRP!M29: was -207
$ws.Cells.Item(29, 14).ClearContents()  # CRP!N29: was -1285.6667

$ws.Cells.Item(31, 8).Value = 14193.5  # CRP!H31: was 12569.028
$ws.Cells.Item(31, 10).Value = 57940.168  # CRP!J31: was 32886.547
$ws.Cells.Item(31, 12).Value = 57940.168  # CRP!L31: was 32886.547
$ws.Cells.Item(31, 14).Value = -58530.168  # CRP!N31: was -33476.547

$ws.Cells.Item(34, 8).Value = 14193.5  # CRP!H34: was 12569.028
$ws.Cells.Item(34, 10).Value = 57940.168  # CRP!J34: was 32886.547
$ws.Cells.Item(34, 12).Value = 57940.168  # CRP!L34: was 32886.547
$ws.Cells.Item(34, 14).Value = -58344.168  # CRP!N34: was -33290.547

$ws.Cells.Item(105, 8).Value = 1445.5  # CRP!H105: was 0
$ws.Cells.Item(105, 9).Value = 885.5  # CRP!I105: was 0
$ws.Cells.Item(105, 10).Value = 2005.5  # CRP!J105: was 0
$ws.Cells.Item(105, 11).Value = 885.5  # CRP!K105: was 0
$ws.Cells.Item(105, 12).Value = 2005.5  # CRP!L105: was 0
$ws.Cells.Item(105, 13).Value = 861.5  # CRP!M105: was None
$ws.Cells.Item(105, 14).Value = -5499.5  # CRP!N105: was None

$ws.Cells.Item(122, 8).Value = 125923.25  # CRP!H122: was 77998.766
$ws.Cells.Item(122, 9).Value = 125923.25  # CRP!I122: was 125944
$ws.Cells.Item(122, 10).Value = 0  # CRP!J122: was 1286.4
$ws.Cells.Item(122, 11).Value = 377769.75  # CRP!K122: was 377832
$ws.Cells.Item(122, 12).Value = 0  # CRP!L122: was 3859.2
$ws.Cells.Item(122, 13).Value = -375319.75  # CRP!M122: was -375382
$ws.Cells.Item(122, 14).ClearContents()  # CRP!N122: was -8759.200000000001

$ws.Cells.Item(132, 8).Value = 4397.222  # CRP!H132: was 4560.375
$ws.Cells.Item(132, 9).Value = 4328.5625  # CRP!I132: was 4505.2144
$ws.Cells.Item(132, 11).Value = 12985.6875  # CRP!K132: was 13515.6432
$ws.Cells.Item(132, 13).Value = -10455.6875  # CRP!M132: was -10985.6432

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 307.14285  # CUL!H10: was 221.875
$ws.Cells.Item(10, 9).Value = 275  # CUL!I10: was 182.14285
$ws.Cells.Item(10, 11).Value = 825  # CUL!K10: was 546.4285500000001
$ws.Cells.Item(10, 13).Value = -686  # CUL!M10: was -407.4285500000001

$ws.Cells.Item(23, 8).Value = 117.84615  # CUL!H23: was 112.92857
$ws.Cells.Item(23, 10).Value = 107  # CUL!J23: was 98.71429000000001
$ws.Cells.Item(23, 12).Value = 321  # CUL!L23: was 296.14287
$ws.Cells.Item(23, 14).Value = -791  # CUL!N23: was -766.14287

$ws.Cells.Item(41, 8).Value = 120.25  # CUL!H41: was 128.95653
$ws.Cells.Item(41, 9).Value = 90.94444  # CUL!I41: was 96
$ws.Cells.Item(41, 10).Value = 173  # CUL!J41: was 204.28572
$ws.Cells.Item(41, 11).Value = 272.83332  # CUL!K41: was 288
$ws.Cells.Item(41, 12).Value = 519  # CUL!L41: was 612.85716
$ws.Cells.Item(41, 13).Value = 65.16667999999999  # CUL!M41: was 50
$ws.Cells.Item(41, 14).Value = -1195  # CUL!N41: was -1288.85716

$ws.Cells.Item(56, 8).Value = 9058.105  # CUL!H56: was 8991.549999999999
$ws.Cells.Item(56, 9).Value = 9058.105  # CUL!I56: was 8991.549999999999
$ws.Cells.Item(56, 11).Value = 9058.105  # CUL!K56: was 8991.549999999999
$ws.Cells.Item(56, 13).Value = -8528.105  # CUL!M56: was -8461.549999999999

$ws.Cells.Item(131, 8).Value = 2210.2  # CUL!H131: was 1821.1578
$ws.Cells.Item(131, 9).Value = 3015  # CUL!I131: was 1246.5
$ws.Cells.Item(131, 11).Value = 9045  # CUL!K131: was 3739.5
$ws.Cells.Item(131, 13).Value = -4005  # CUL!M131: was 1300.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9531.091  # GSM!H70: was 9772.666999999999
$ws.Cells.Item(70, 9).Value = 9605.375  # GSM!I70: was 9992.5
$ws.Cells.Item(70, 11).Value = 9605.375  # GSM!K70: was 9992.5
$ws.Cells.Item(70, 13).Value = -9335.375  # GSM!M70: was -9722.5

$ws.Cells.Item(73, 8).Value = 9531.091  # GSM!H73: was 9772.666999999999
$ws.Cells.Item(73, 9).Value = 9605.375  # GSM!I73: was 9992.5
$ws.Cells.Item(73, 11).Value = 9605.375  # GSM!K73: was 9992.5
$ws.Cells.Item(73, 13).Value = -8669.375  # GSM!M73: was -9056.5

$ws.Cells.Item(132, 8).Value = 11788.5  # GSM!H132: was 11746.2

$ws.Cells.Item(141, 8).Value = 66714.5  # GSM!H141: was 61143
$ws.Cells.Item(141, 10).Value = 66714.5  # GSM!J141: was 61143
$ws.Cells.Item(141, 12).Value = 66714.5  # GSM!L141: was 61143
$ws.Cells.Item(141, 14).Value = -77074.5  # GSM!N141: was -71503

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1047.75  # LTW!H22: was 1026.6666
$ws.Cells.Item(22, 10).Value = 1111  # LTW!J22: was 0
$ws.Cells.Item(22, 12).Value = 1111  # LTW!L22: was 0
$ws.Cells.Item(22, 14).Value = -1701  # LTW!N22: was None

$ws.Cells.Item(27, 8).Value = 1047.75  # LTW!H27: was 1026.6666
$ws.Cells.Item(27, 10).Value = 1111  # LTW!J27: was 0
$ws.Cells.Item(27, 12).Value = 1111  # LTW!L27: was 0
$ws.Cells.Item(27, 14).Value = -1325  # LTW!N27: was None

$ws.Cells.Item(46, 8).Value = 25861.666  # LTW!H46: was 25881.223
$ws.Cells.Item(46, 9).Value = 72012.164  # LTW!I46: was 72070.336
$ws.Cells.Item(46, 10).Value = 2786.4167  # LTW!J46: was 2786.6667
$ws.Cells.Item(46, 11).Value = 72012.164  # LTW!K46: was 72070.336
$ws.Cells.Item(46, 12).Value = 2786.4167  # LTW!L46: was 2786.6667
$ws.Cells.Item(46, 13).Value = -71824.164  # LTW!M46: was -71882.336
$ws.Cells.Item(46, 14).Value = -3162.4167  # LTW!N46: was -3162.6667

$ws.Cells.Item(132, 8).Value = 3787.1853  # LTW!H132: was 3862.8
$ws.Cells.Item(132, 9).Value = 3655.0454  # LTW!I132: was 3736.35
$ws.Cells.Item(132, 11).Value = 10965.1362  # LTW!K132: was 11209.05
$ws.Cells.Item(132, 13).Value = -8435.136200000001  # LTW!M132: was -8679.049999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 34999.5  # WVR!H4: was 21399.75
$ws.Cells.Item(4, 9).Value = 40000  # WVR!I4: was 18533.334
$ws.Cells.Item(4, 11).Value = 40000  # WVR!K4: was 18533.334
$ws.Cells.Item(4, 13).Value = -39887  # WVR!M4: was -18420.334

$ws.Cells.Item(140, 8).Value = 83819.8  # WVR!H140: was 84775
$ws.Cells.Item(140, 10).Value = 83499.75  # WVR!J140: was 84666.664
$ws.Cells.Item(140, 12).Value = 83499.75  # WVR!L140: was 84666.664
$ws.Cells.Item(140, 14).Value = -93859.75  # WVR!N140: was -95026.664

$ws.Cells.Item(141, 8).Value = 57853.715  # WVR!H141: was 59218.832
$ws.Cells.Item(141, 10).Value = 60721  # WVR!J141: was 62932.6
$ws.Cells.Item(141, 12).Value = 60721  # WVR!L141: was 62932.6
$ws.Cells.Item(141, 14).Value = -71081  # WVR!N141: was -73292.60000000001
